# Update the footer "Update automatically" date field from 13.01.2022 to
# 17.01.2022 across the slide master, every slide layout, the handout
# master and the notes master, and fix a typo in the Coin class
# description on slide 6 ("обломки ключа" -> "монеты").

$p = $ppt.ActivePresentation

$oldDate = "13.01.2022"
$newDate = "17.01.2022"

function Update-DateShapes($container) {
    $shapes = $container.Shapes
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# Slide master footer date placeholder.
Update-DateShapes $p.SlideMaster

# Every slide layout's footer date placeholder.
$layouts = $p.SlideMaster.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    Update-DateShapes $layouts.Item($L)
}

# Handout master and notes master: direct shape-text edits aren't
# applied for these two containers, but the Header/Footer "date and
# time" object is, so go through HeadersFooters.DateAndTime instead.
$p.HandoutMaster.HeadersFooters.DateAndTime.Text = $newDate
$p.NotesMaster.HeadersFooters.DateAndTime.Text = $newDate

# Slide 6: "class Coin" description text tweak.
$slide6 = $p.Slides.Item(6)
$shapes6 = $slide6.Shapes
for ($i = 1; $i -le $shapes6.Count; $i++) {
    $sh = $shapes6.Item($i)
    if ($sh.HasTextFrame) {
        $t = $sh.TextFrame.TextRange.Text
        if ($t -eq "В этом классе реализованы очки(обломки ключа), который надо собрать на локации чтобы открыть выход и пройти игру.") {
            $sh.TextFrame.TextRange.Text = "В этом классе реализованы очки(монеты), который надо собрать на локации чтобы открыть выход и пройти игру."
        }
    }
}
